$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.6292110000000001
$ws.Range("H2").Value = 1.887633
$ws.Range("I2").Value = 0.01078649253029594
$ws.Range("J2").Value = 0.01078649253029594
$ws.Range("M2").Value = 0.353079
$ws.Range("N2").Value = 1.059237
$ws.Range("O2").Value = 0.01390801122570493
$ws.Range("P2").Value = 0.01390801122570493
$ws.Range("Q2").Value = 0.222161190669
$ws.Range("R2").Value = 1.999450716021
$ws.Range("S2").Value = 0.0001500186591973383
$ws.Range("T2").Value = 0.0001500186591973383

# Row 3
$ws.Range("G3").Value = 0.6292110000000001
$ws.Range("H3").Value = 1.887633
$ws.Range("I3").Value = 0.01078649253029594
$ws.Range("J3").Value = 0.01078649253029594
$ws.Range("M3").Value = 24.359095
$ws.Range("N3").Value = 73.077285
$ws.Range("O3").Value = 0.9595205795530543
$ws.Range("P3").Value = 0.9595205795530543
$ws.Range("Q3").Value = 15.327010524045
$ws.Range("R3").Value = 137.943094716405
$ws.Range("S3").Value = 0.01034986156401425
$ws.Range("T3").Value = 0.01034986156401425

# Row 4
$ws.Range("G4").Value = 0.6292110000000001
$ws.Range("H4").Value = 1.887633
$ws.Range("I4").Value = 0.01078649253029594
$ws.Range("J4").Value = 0.01078649253029594
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.6745613333333332
$ws.Range("N4").Value = 2.023684
$ws.Range("O4").Value = 0.02657140922124081
$ws.Range("P4").Value = 0.02657140922124081
$ws.Range("Q4").Value = 0.424441411108
$ws.Range("R4").Value = 3.819972699972
$ws.Range("S4").Value = 0.0002866123070843506
$ws.Range("T4").Value = 0.0002866123070843506

# Row 5
$ws.Range("I5").Value = 0.4124494462002061
$ws.Range("J5").Value = 0.4124494462002061
$ws.Range("M5").Value = 0.353079
$ws.Range("N5").Value = 1.059237
$ws.Range("O5").Value = 0.01390801122570493
$ws.Range("P5").Value = 0.01390801122570493
$ws.Range("Q5").Value = 8.494907848982999
$ws.Range("R5").Value = 76.454170640847
$ws.Range("S5").Value = 0.005736351527788248
$ws.Range("T5").Value = 0.005736351527788248

# Row 6
$ws.Range("I6").Value = 0.4124494462002061
$ws.Range("J6").Value = 0.4124494462002061
$ws.Range("M6").Value = 24.359095
$ws.Range("N6").Value = 73.077285
$ws.Range("O6").Value = 0.9595205795530543
$ws.Range("P6").Value = 0.9595205795530543
$ws.Range("Q6").Value = 586.0678978631483
$ws.Range("R6").Value = 5274.611080768334
$ws.Range("S6").Value = 0.3957537316543581
$ws.Range("T6").Value = 0.395753731654358

# Row 7
$ws.Range("I7").Value = 0.4124494462002061
$ws.Range("J7").Value = 0.4124494462002061
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.6745613333333332
$ws.Range("N7").Value = 2.023684
$ws.Range("O7").Value = 0.02657140922124081
$ws.Range("P7").Value = 0.02657140922124081
$ws.Range("Q7").Value = 16.22961536980044
$ws.Range("R7").Value = 146.066538328204
$ws.Range("S7").Value = 0.01095936301805982
$ws.Range("T7").Value = 0.01095936301805982

# Row 8
$ws.Range("G8").Value = 16.061552
$ws.Range("H8").Value = 48.184656
$ws.Range("I8").Value = 0.2753413571488099
$ws.Range("J8").Value = 0.2753413571488098
$ws.Range("M8").Value = 0.353079
$ws.Range("N8").Value = 1.059237
$ws.Range("O8").Value = 0.01390801122570493
$ws.Range("P8").Value = 0.01390801122570493
$ws.Range("Q8").Value = 5.670996718608
$ws.Range("R8").Value = 51.03897046747201
$ws.Range("S8").Value = 0.003829450686126478
$ws.Range("T8").Value = 0.003829450686126477

# Row 9
$ws.Range("G9").Value = 16.061552
$ws.Range("H9").Value = 48.184656
$ws.Range("I9").Value = 0.2753413571488099
$ws.Range("J9").Value = 0.2753413571488098
$ws.Range("M9").Value = 24.359095
$ws.Range("N9").Value = 73.077285
$ws.Range("O9").Value = 0.9595205795530543
$ws.Range("P9").Value = 0.9595205795530543
$ws.Range("Q9").Value = 391.2448710154401
$ws.Range("R9").Value = 3521.203839138961
$ws.Range("S9").Value = 0.2641956985863506
$ws.Range("T9").Value = 0.2641956985863505

# Row 10
$ws.Range("G10").Value = 16.061552
$ws.Range("H10").Value = 48.184656
$ws.Range("I10").Value = 0.2753413571488099
$ws.Range("J10").Value = 0.2753413571488098
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.6745613333333332
$ws.Range("N10").Value = 2.023684
$ws.Range("O10").Value = 0.02657140922124081
$ws.Range("P10").Value = 0.02657140922124081
$ws.Range("Q10").Value = 10.83450193252267
$ws.Range("R10").Value = 97.510517392704
$ws.Range("S10").Value = 0.007316207876332846
$ws.Range("T10").Value = 0.007316207876332845

# Row 11
$ws.Range("G11").Value = 1.855556
$ws.Range("H11").Value = 5.566668
$ws.Range("I11").Value = 0.03180958523221274
$ws.Range("J11").Value = 0.03180958523221274
$ws.Range("M11").Value = 0.353079
$ws.Range("N11").Value = 1.059237
$ws.Range("O11").Value = 0.01390801122570493
$ws.Range("P11").Value = 0.01390801122570493
$ws.Range("Q11").Value = 0.655157856924
$ws.Range("R11").Value = 5.896420712316
$ws.Range("S11").Value = 0.0004424080684946325
$ws.Range("T11").Value = 0.0004424080684946325

# Row 12
$ws.Range("G12").Value = 1.855556
$ws.Range("H12").Value = 5.566668
$ws.Range("I12").Value = 0.03180958523221274
$ws.Range("J12").Value = 0.03180958523221274
$ws.Range("M12").Value = 24.359095
$ws.Range("N12").Value = 73.077285
$ws.Range("O12").Value = 0.9595205795530543
$ws.Range("P12").Value = 0.9595205795530543
$ws.Range("Q12").Value = 45.19966488182
$ws.Range("R12").Value = 406.79698393638
$ws.Range("S12").Value = 0.03052195165735504
$ws.Range("T12").Value = 0.03052195165735504

# Row 13
$ws.Range("G13").Value = 1.855556
$ws.Range("H13").Value = 5.566668
$ws.Range("I13").Value = 0.03180958523221274
$ws.Range("J13").Value = 0.03180958523221274
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.6745613333333332
$ws.Range("N13").Value = 2.023684
$ws.Range("O13").Value = 0.02657140922124081
$ws.Range("P13").Value = 0.02657140922124081
$ws.Range("Q13").Value = 1.251686329434667
$ws.Range("R13").Value = 11.265176964912
$ws.Range("S13").Value = 0.000845225506363063
$ws.Range("T13").Value = 0.000845225506363063

# Row 14
$ws.Range("G14").Value = 15.72740533333333
$ws.Range("H14").Value = 47.182216
$ws.Range("I14").Value = 0.2696131188884753
$ws.Range("J14").Value = 0.2696131188884753
$ws.Range("M14").Value = 0.353079
$ws.Range("N14").Value = 1.059237
$ws.Range("O14").Value = 0.01390801122570493
$ws.Range("P14").Value = 0.01390801122570493
$ws.Range("Q14").Value = 5.553016547687999
$ws.Range("R14").Value = 49.977148929192
$ws.Range("S14").Value = 0.003749782284098232
$ws.Range("T14").Value = 0.003749782284098233

# Row 15
$ws.Range("G15").Value = 15.72740533333333
$ws.Range("H15").Value = 47.182216
$ws.Range("I15").Value = 0.2696131188884753
$ws.Range("J15").Value = 0.2696131188884753
$ws.Range("M15").Value = 24.359095
$ws.Range("N15").Value = 73.077285
$ws.Range("O15").Value = 0.9595205795530543
$ws.Range("P15").Value = 0.9595205795530543
$ws.Range("Q15").Value = 383.1053606181733
$ws.Range("R15").Value = 3447.94824556356
$ws.Range("S15").Value = 0.2586993360909764
$ws.Range("T15").Value = 0.2586993360909764

# Row 16
$ws.Range("G16").Value = 15.72740533333333
$ws.Range("H16").Value = 47.182216
$ws.Range("I16").Value = 0.2696131188884753
$ws.Range("J16").Value = 0.2696131188884753
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 0.6745613333333332
$ws.Range("N16").Value = 2.023684
$ws.Range("O16").Value = 0.02657140922124081
$ws.Range("P16").Value = 0.02657140922124081
$ws.Range("Q16").Value = 10.60909951152711
$ws.Range("R16").Value = 95.48189560374398
$ws.Range("S16").Value = 0.007164000513400728
$ws.Range("T16").Value = 0.007164000513400728
